# Applies the crypto price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.094.40"
$ws.Range("E2").Value = "  -3.55%  "

$ws.Range("D3").Value = "'1.651.14"
$ws.Range("E3").Value = "  -5.39%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'237.61"
$ws.Range("E5").Value = "  -4.78%  "

$ws.Range("E6").Value = "  +0.04%  "

$ws.Range("D7").Value = "'0.4794"
$ws.Range("E7").Value = "  -7.20%  "

$ws.Range("D8").Value = "'0.2610"
$ws.Range("E8").Value = "  -5.56%  "

$ws.Range("D9").Value = "'0.06005"
$ws.Range("E9").Value = "  -3.11%  "

$ws.Range("D10").Value = "'0.07185"
$ws.Range("E10").Value = "  -0.52%  "

$ws.Range("D11").Value = "'1.657.48"
$ws.Range("E11").Value = "  -4.84%  "

$ws.Range("D12").Value = "'14.85"
$ws.Range("E12").Value = "  -2.58%  "

$ws.Range("D13").Value = "'0.6236"
$ws.Range("E13").Value = "  -4.16%  "

$ws.Range("D14").Value = "'4.605"
$ws.Range("E14").Value = "  -0.70%  "

$ws.Range("D15").Value = "'73.42"
$ws.Range("E15").Value = "  -5.85%  "

$ws.Range("D16").Value = "'0.9996"
$ws.Range("E16").Value = "  -0.05%  "

$ws.Range("D17").Value = "'1.000"
$ws.Range("E17").Value = "  +0.08%  "

$ws.Range("D18").Value = "'25.079.56"
$ws.Range("E18").Value = "  -3.74%  "

$ws.Range("D19").Value = "'11.45"
$ws.Range("E19").Value = "  -3.53%  "

$ws.Range("D20").Value = "'0.000006600"
$ws.Range("E20").Value = "  -3.19%  "

$ws.Range("D21").Value = "'4.470"
$ws.Range("E21").Value = "  +4.04%  "

$ws.Range("D22").Value = "'1.863.85"
$ws.Range("E22").Value = "  -5.25%  "

$ws.Range("D23").Value = "'8.607"
$ws.Range("E23").Value = "  -0.98%  "

$ws.Range("D24").Value = "'5.289"
$ws.Range("E24").Value = "  -1.50%  "

$ws.Range("D25").Value = "'133.08"
$ws.Range("E25").Value = "  -1.93%  "

$ws.Range("D26").Value = "'14.94"
$ws.Range("E26").Value = "  -2.37%  "

$ws.Range("D27").Value = "'1.394"
$ws.Range("E27").Value = "  -7.51%  "

$ws.Range("D28").Value = "'103.41"
$ws.Range("E28").Value = "  -2.28%  "

$ws.Range("D29").Value = "'1.680"
$ws.Range("E29").Value = "  -5.71%  "

$ws.Range("D30").Value = "'3.772"
$ws.Range("E30").Value = "  -4.82%  "

$ws.Range("D31").Value = "'0.07918"
$ws.Range("E31").Value = "  -4.41%  "

$ws.Range("D32").Value = "'3.568"
$ws.Range("E32").Value = "  -3.14%  "

$ws.Range("D33").Value = "'0.04594"
$ws.Range("E33").Value = "  -1.71%  "

$ws.Range("D34").Value = "'2.594"
$ws.Range("E34").Value = "  -2.33%  "

$ws.Range("D35").Value = "'0.9450"
$ws.Range("E35").Value = "  -5.77%  "

$ws.Range("D36").Value = "'0.5775"
$ws.Range("E36").Value = "  -7.41%  "

$ws.Range("D37").Value = "'2.613"
$ws.Range("E37").Value = "  -4.71%  "

$ws.Range("D38").Value = "'0.01554"
$ws.Range("E38").Value = "  -3.14%  "

$ws.Range("B39").Value = "PaxDollar"
$ws.Range("C39").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D39").Value = "'1.000"
$ws.Range("E39").Value = "  +0.10%  "

$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "'0.8340"
$ws.Range("E40").Value = "  +10.06%  "

$ws.Range("D41").Value = "'1.828"
$ws.Range("E41").Value = "  -5.69%  "

$ws.Range("D42").Value = "'99.07"
$ws.Range("E42").Value = "  -1.46%  "

$ws.Range("E43").Value = "  -4.24%  "

$ws.Range("D44").Value = "'4.815"
$ws.Range("E44").Value = "  -3.89%  "

$ws.Range("D45").Value = "'0.1138"
$ws.Range("E45").Value = "  +0.31%  "

$ws.Range("D46").Value = "'6.108"
$ws.Range("E46").Value = "  -3.87%  "

$ws.Range("D47").Value = "'0.05185"
$ws.Range("E47").Value = "  -0.88%  "

$ws.Range("D48").Value = "'29.82"
$ws.Range("E48").Value = "  -2.93%  "

$ws.Range("D49").Value = "'51.17"
$ws.Range("E49").Value = "  -7.86%  "

$ws.Range("D50").Value = "'1.002"
$ws.Range("E50").Value = "  -0.03%  "

$ws.Range("D51").Value = "'0.3339"
$ws.Range("E51").Value = "  -3.06%  "
